# errors on register page fixed after rec system release
# Update the promo-code related label/help text on the "register" sheet
# and move the active selection to C65 (matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63 (label_promo_code): new English/Dutch label text
$ws.Range("B63").Value = "Did a friend give you a promocode? (OPTIONAL) "
$ws.Range("C63").Value = "Heeft een vriend u een promocode gegeven? (OPTIONEEL)"

# Row 65 (valid_feedback_promo_code): shortened Dutch helper text
$ws.Range("C65").Value = "Deze code wordt gebruikt om degene die u heeft doorverwezen te belonen. Door de beloning zou hij/zij kunnen zien of u wel of niet succesvol heeft deelgenomen aan dit onderzoek."

# Update the active selection to C65
$ws.Range("C65").Select()
